# Update DAIRY - MODIFIED LAND (dry: I, irr: AQ) and DAIRY - NATURAL LAND (dry: J, irr: AR)
# productivity/FLC values for rows 6-43 across all three scenario sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Business As Usual")
$ws.Range("I6").Value = 0.9959492779147587
$ws.Range("J6").Value = 0.9959492779147587
$ws.Range("AQ6").Value = 0.9959492779147587
$ws.Range("AR6").Value = 0.9959492779147587
$ws.Range("I7").Value = 0.9918985558295175
$ws.Range("J7").Value = 0.9918985558295175
$ws.Range("AQ7").Value = 0.9918985558295175
$ws.Range("AR7").Value = 0.9918985558295175
$ws.Range("I8").Value = 0.9878478337442762
$ws.Range("J8").Value = 0.9878478337442762
$ws.Range("AQ8").Value = 0.9878478337442762
$ws.Range("AR8").Value = 0.9878478337442762
$ws.Range("I9").Value = 0.9837971116590348
$ws.Range("J9").Value = 0.9837971116590348
$ws.Range("AQ9").Value = 0.9837971116590348
$ws.Range("AR9").Value = 0.9837971116590348
$ws.Range("I10").Value = 0.9797463895737936
$ws.Range("J10").Value = 0.9797463895737936
$ws.Range("AQ10").Value = 0.9797463895737936
$ws.Range("AR10").Value = 0.9797463895737936
$ws.Range("I11").Value = 0.9756956674885523
$ws.Range("J11").Value = 0.9756956674885523
$ws.Range("AQ11").Value = 0.9756956674885523
$ws.Range("AR11").Value = 0.9756956674885523
$ws.Range("I12").Value = 0.971644945403311
$ws.Range("J12").Value = 0.971644945403311
$ws.Range("AQ12").Value = 0.971644945403311
$ws.Range("AR12").Value = 0.971644945403311
$ws.Range("I13").Value = 0.9675942233180698
$ws.Range("J13").Value = 0.9675942233180698
$ws.Range("AQ13").Value = 0.9675942233180698
$ws.Range("AR13").Value = 0.9675942233180698
$ws.Range("I14").Value = 0.9635435012328284
$ws.Range("J14").Value = 0.9635435012328284
$ws.Range("AQ14").Value = 0.9635435012328284
$ws.Range("AR14").Value = 0.9635435012328284
$ws.Range("I15").Value = 0.9594927791475871
$ws.Range("J15").Value = 0.9594927791475871
$ws.Range("AQ15").Value = 0.9594927791475871
$ws.Range("AR15").Value = 0.9594927791475871
$ws.Range("I16").Value = 0.9554420570623459
$ws.Range("J16").Value = 0.9554420570623459
$ws.Range("AQ16").Value = 0.9554420570623459
$ws.Range("AR16").Value = 0.9554420570623459
$ws.Range("I17").Value = 0.9513913349771046
$ws.Range("J17").Value = 0.9513913349771046
$ws.Range("AQ17").Value = 0.9513913349771046
$ws.Range("AR17").Value = 0.9513913349771046
$ws.Range("I18").Value = 0.9473406128918633
$ws.Range("J18").Value = 0.9473406128918633
$ws.Range("AQ18").Value = 0.9473406128918633
$ws.Range("AR18").Value = 0.9473406128918633
$ws.Range("I19").Value = 0.9432898908066221
$ws.Range("J19").Value = 0.9432898908066221
$ws.Range("AQ19").Value = 0.9432898908066221
$ws.Range("AR19").Value = 0.9432898908066221
$ws.Range("I20").Value = 0.9392391687213808
$ws.Range("J20").Value = 0.9392391687213808
$ws.Range("AQ20").Value = 0.9392391687213808
$ws.Range("AR20").Value = 0.9392391687213808
$ws.Range("I21").Value = 0.9351884466361394
$ws.Range("J21").Value = 0.9351884466361394
$ws.Range("AQ21").Value = 0.9351884466361394
$ws.Range("AR21").Value = 0.9351884466361394
$ws.Range("I22").Value = 0.9311377245508982
$ws.Range("J22").Value = 0.9311377245508982
$ws.Range("AQ22").Value = 0.9311377245508982
$ws.Range("AR22").Value = 0.9311377245508982
$ws.Range("I23").Value = 0.9270870024656569
$ws.Range("J23").Value = 0.9270870024656569
$ws.Range("AQ23").Value = 0.9270870024656569
$ws.Range("AR23").Value = 0.9270870024656569
$ws.Range("I24").Value = 0.9231067277210285
$ws.Range("J24").Value = 0.9231067277210285
$ws.Range("AQ24").Value = 0.9231067277210285
$ws.Range("AR24").Value = 0.9231067277210285
$ws.Range("I25").Value = 0.9191264529764002
$ws.Range("J25").Value = 0.9191264529764002
$ws.Range("AQ25").Value = 0.9191264529764002
$ws.Range("AR25").Value = 0.9191264529764002
$ws.Range("I26").Value = 0.9151461782317717
$ws.Range("J26").Value = 0.9151461782317717
$ws.Range("AQ26").Value = 0.9151461782317717
$ws.Range("AR26").Value = 0.9151461782317717
$ws.Range("I27").Value = 0.9111659034871433
$ws.Range("J27").Value = 0.9111659034871433
$ws.Range("AQ27").Value = 0.9111659034871433
$ws.Range("AR27").Value = 0.9111659034871433
$ws.Range("I28").Value = 0.907185628742515
$ws.Range("J28").Value = 0.907185628742515
$ws.Range("AQ28").Value = 0.907185628742515
$ws.Range("AR28").Value = 0.907185628742515
$ws.Range("I29").Value = 0.9037689327227897
$ws.Range("J29").Value = 0.9037689327227897
$ws.Range("AQ29").Value = 0.9037689327227897
$ws.Range("AR29").Value = 0.9037689327227897
$ws.Range("I30").Value = 0.9003522367030644
$ws.Range("J30").Value = 0.9003522367030644
$ws.Range("AQ30").Value = 0.9003522367030644
$ws.Range("AR30").Value = 0.9003522367030644
$ws.Range("I31").Value = 0.8969355406833392
$ws.Range("J31").Value = 0.8969355406833392
$ws.Range("AQ31").Value = 0.8969355406833392
$ws.Range("AR31").Value = 0.8969355406833392
$ws.Range("I32").Value = 0.893518844663614
$ws.Range("J32").Value = 0.893518844663614
$ws.Range("AQ32").Value = 0.893518844663614
$ws.Range("AR32").Value = 0.893518844663614
$ws.Range("I33").Value = 0.8901021486438887
$ws.Range("J33").Value = 0.8901021486438887
$ws.Range("AQ33").Value = 0.8901021486438887
$ws.Range("AR33").Value = 0.8901021486438887
$ws.Range("I34").Value = 0.8873899260302923
$ws.Range("J34").Value = 0.8873899260302923
$ws.Range("AQ34").Value = 0.8873899260302923
$ws.Range("AR34").Value = 0.8873899260302923
$ws.Range("I35").Value = 0.884677703416696
$ws.Range("J35").Value = 0.884677703416696
$ws.Range("AQ35").Value = 0.884677703416696
$ws.Range("AR35").Value = 0.884677703416696
$ws.Range("I36").Value = 0.8819654808030997
$ws.Range("J36").Value = 0.8819654808030997
$ws.Range("AQ36").Value = 0.8819654808030997
$ws.Range("AR36").Value = 0.8819654808030997
$ws.Range("I37").Value = 0.8792532581895034
$ws.Range("J37").Value = 0.8792532581895034
$ws.Range("AQ37").Value = 0.8792532581895034
$ws.Range("AR37").Value = 0.8792532581895034
$ws.Range("I38").Value = 0.876541035575907
$ws.Range("J38").Value = 0.876541035575907
$ws.Range("AQ38").Value = 0.876541035575907
$ws.Range("AR38").Value = 0.876541035575907
$ws.Range("I39").Value = 0.8738288129623106
$ws.Range("J39").Value = 0.8738288129623106
$ws.Range("AQ39").Value = 0.8738288129623106
$ws.Range("AR39").Value = 0.8738288129623106
$ws.Range("I40").Value = 0.8711165903487144
$ws.Range("J40").Value = 0.8711165903487144
$ws.Range("AQ40").Value = 0.8711165903487144
$ws.Range("AR40").Value = 0.8711165903487144
$ws.Range("I41").Value = 0.868404367735118
$ws.Range("J41").Value = 0.868404367735118
$ws.Range("AQ41").Value = 0.868404367735118
$ws.Range("AR41").Value = 0.868404367735118
$ws.Range("I42").Value = 0.8656921451215217
$ws.Range("J42").Value = 0.8656921451215217
$ws.Range("AQ42").Value = 0.8656921451215217
$ws.Range("AR42").Value = 0.8656921451215217
$ws.Range("I43").Value = 0.8629799225079253
$ws.Range("J43").Value = 0.8629799225079253
$ws.Range("AQ43").Value = 0.8629799225079253
$ws.Range("AR43").Value = 0.8629799225079253

$ws = $wb.Worksheets.Item("Stratified Societies")
$ws.Range("I6").Value = 0.9927106571171382
$ws.Range("J6").Value = 0.9927106571171382
$ws.Range("AQ6").Value = 0.9927106571171382
$ws.Range("AR6").Value = 0.9927106571171382
$ws.Range("I7").Value = 0.9854213142342766
$ws.Range("J7").Value = 0.9854213142342766
$ws.Range("AQ7").Value = 0.9854213142342766
$ws.Range("AR7").Value = 0.9854213142342766
$ws.Range("I8").Value = 0.9781319713514148
$ws.Range("J8").Value = 0.9781319713514148
$ws.Range("AQ8").Value = 0.9781319713514148
$ws.Range("AR8").Value = 0.9781319713514148
$ws.Range("I9").Value = 0.9708426284685531
$ws.Range("J9").Value = 0.9708426284685531
$ws.Range("AQ9").Value = 0.9708426284685531
$ws.Range("AR9").Value = 0.9708426284685531
$ws.Range("I10").Value = 0.9635532855856913
$ws.Range("J10").Value = 0.9635532855856913
$ws.Range("AQ10").Value = 0.9635532855856913
$ws.Range("AR10").Value = 0.9635532855856913
$ws.Range("I11").Value = 0.9562639427028297
$ws.Range("J11").Value = 0.9562639427028297
$ws.Range("AQ11").Value = 0.9562639427028297
$ws.Range("AR11").Value = 0.9562639427028297
$ws.Range("I12").Value = 0.9489745998199679
$ws.Range("J12").Value = 0.9489745998199679
$ws.Range("AQ12").Value = 0.9489745998199679
$ws.Range("AR12").Value = 0.9489745998199679
$ws.Range("I13").Value = 0.9416852569371061
$ws.Range("J13").Value = 0.9416852569371061
$ws.Range("AQ13").Value = 0.9416852569371061
$ws.Range("AR13").Value = 0.9416852569371061
$ws.Range("I14").Value = 0.9343959140542444
$ws.Range("J14").Value = 0.9343959140542444
$ws.Range("AQ14").Value = 0.9343959140542444
$ws.Range("AR14").Value = 0.9343959140542444
$ws.Range("I15").Value = 0.9271065711713827
$ws.Range("J15").Value = 0.9271065711713827
$ws.Range("AQ15").Value = 0.9271065711713827
$ws.Range("AR15").Value = 0.9271065711713827
$ws.Range("I16").Value = 0.919817228288521
$ws.Range("J16").Value = 0.919817228288521
$ws.Range("AQ16").Value = 0.919817228288521
$ws.Range("AR16").Value = 0.919817228288521
$ws.Range("I17").Value = 0.9125278854056592
$ws.Range("J17").Value = 0.9125278854056592
$ws.Range("AQ17").Value = 0.9125278854056592
$ws.Range("AR17").Value = 0.9125278854056592
$ws.Range("I18").Value = 0.9052385425227976
$ws.Range("J18").Value = 0.9052385425227976
$ws.Range("AQ18").Value = 0.9052385425227976
$ws.Range("AR18").Value = 0.9052385425227976
$ws.Range("I19").Value = 0.8979491996399358
$ws.Range("J19").Value = 0.8979491996399358
$ws.Range("AQ19").Value = 0.8979491996399358
$ws.Range("AR19").Value = 0.8979491996399358
$ws.Range("I20").Value = 0.890659856757074
$ws.Range("J20").Value = 0.890659856757074
$ws.Range("AQ20").Value = 0.890659856757074
$ws.Range("AR20").Value = 0.890659856757074
$ws.Range("I21").Value = 0.8833705138742123
$ws.Range("J21").Value = 0.8833705138742123
$ws.Range("AQ21").Value = 0.8833705138742123
$ws.Range("AR21").Value = 0.8833705138742123
$ws.Range("I22").Value = 0.8760811709913506
$ws.Range("J22").Value = 0.8760811709913506
$ws.Range("AQ22").Value = 0.8760811709913506
$ws.Range("AR22").Value = 0.8760811709913506
$ws.Range("I23").Value = 0.8687918281084889
$ws.Range("J23").Value = 0.8687918281084889
$ws.Range("AQ23").Value = 0.8687918281084889
$ws.Range("AR23").Value = 0.8687918281084889
$ws.Range("I24").Value = 0.8658330398027474
$ws.Range("J24").Value = 0.8658330398027474
$ws.Range("AQ24").Value = 0.8658330398027474
$ws.Range("AR24").Value = 0.8658330398027474
$ws.Range("I25").Value = 0.8628742514970059
$ws.Range("J25").Value = 0.8628742514970059
$ws.Range("AQ25").Value = 0.8628742514970059
$ws.Range("AR25").Value = 0.8628742514970059
$ws.Range("I26").Value = 0.8599154631912646
$ws.Range("J26").Value = 0.8599154631912646
$ws.Range("AQ26").Value = 0.8599154631912646
$ws.Range("AR26").Value = 0.8599154631912646
$ws.Range("I27").Value = 0.8569566748855231
$ws.Range("J27").Value = 0.8569566748855231
$ws.Range("AQ27").Value = 0.8569566748855231
$ws.Range("AR27").Value = 0.8569566748855231
$ws.Range("I28").Value = 0.8539978865797816
$ws.Range("J28").Value = 0.8539978865797816
$ws.Range("AQ28").Value = 0.8539978865797816
$ws.Range("AR28").Value = 0.8539978865797816
$ws.Range("I29").Value = 0.8535047551954914
$ws.Range("J29").Value = 0.8535047551954914
$ws.Range("AQ29").Value = 0.8535047551954914
$ws.Range("AR29").Value = 0.8535047551954914
$ws.Range("I30").Value = 0.8530116238112011
$ws.Range("J30").Value = 0.8530116238112011
$ws.Range("AQ30").Value = 0.8530116238112011
$ws.Range("AR30").Value = 0.8530116238112011
$ws.Range("I31").Value = 0.8525184924269109
$ws.Range("J31").Value = 0.8525184924269109
$ws.Range("AQ31").Value = 0.8525184924269109
$ws.Range("AR31").Value = 0.8525184924269109
$ws.Range("I32").Value = 0.8520253610426206
$ws.Range("J32").Value = 0.8520253610426206
$ws.Range("AQ32").Value = 0.8520253610426206
$ws.Range("AR32").Value = 0.8520253610426206
$ws.Range("I33").Value = 0.8515322296583304
$ws.Range("J33").Value = 0.8515322296583304
$ws.Range("AQ33").Value = 0.8515322296583304
$ws.Range("AR33").Value = 0.8515322296583304
$ws.Range("I34").Value = 0.8526065516026771
$ws.Range("J34").Value = 0.8526065516026771
$ws.Range("AQ34").Value = 0.8526065516026771
$ws.Range("AR34").Value = 0.8526065516026771
$ws.Range("I35").Value = 0.8536808735470236
$ws.Range("J35").Value = 0.8536808735470236
$ws.Range("AQ35").Value = 0.8536808735470236
$ws.Range("AR35").Value = 0.8536808735470236
$ws.Range("I36").Value = 0.8547551954913702
$ws.Range("J36").Value = 0.8547551954913702
$ws.Range("AQ36").Value = 0.8547551954913702
$ws.Range("AR36").Value = 0.8547551954913702
$ws.Range("I37").Value = 0.8558295174357168
$ws.Range("J37").Value = 0.8558295174357168
$ws.Range("AQ37").Value = 0.8558295174357168
$ws.Range("AR37").Value = 0.8558295174357168
$ws.Range("I38").Value = 0.8569038393800634
$ws.Range("J38").Value = 0.8569038393800634
$ws.Range("AQ38").Value = 0.8569038393800634
$ws.Range("AR38").Value = 0.8569038393800634
$ws.Range("I39").Value = 0.8579781613244101
$ws.Range("J39").Value = 0.8579781613244101
$ws.Range("AQ39").Value = 0.8579781613244101
$ws.Range("AR39").Value = 0.8579781613244101
$ws.Range("I40").Value = 0.8590524832687566
$ws.Range("J40").Value = 0.8590524832687566
$ws.Range("AQ40").Value = 0.8590524832687566
$ws.Range("AR40").Value = 0.8590524832687566
$ws.Range("I41").Value = 0.8601268052131033
$ws.Range("J41").Value = 0.8601268052131033
$ws.Range("AQ41").Value = 0.8601268052131033
$ws.Range("AR41").Value = 0.8601268052131033
$ws.Range("I42").Value = 0.8612011271574498
$ws.Range("J42").Value = 0.8612011271574498
$ws.Range("AQ42").Value = 0.8612011271574498
$ws.Range("AR42").Value = 0.8612011271574498
$ws.Range("I43").Value = 0.8622754491017964
$ws.Range("J43").Value = 0.8622754491017964
$ws.Range("AQ43").Value = 0.8622754491017964
$ws.Range("AR43").Value = 0.8622754491017964

$ws = $wb.Worksheets.Item("Toward Sustainability")
$ws.Range("I6").Value = 0.9974756369613714
$ws.Range("J6").Value = 0.9974756369613714
$ws.Range("AQ6").Value = 0.9974756369613714
$ws.Range("AR6").Value = 0.9974756369613714
$ws.Range("I7").Value = 0.9949512739227427
$ws.Range("J7").Value = 0.9949512739227427
$ws.Range("AQ7").Value = 0.9949512739227427
$ws.Range("AR7").Value = 0.9949512739227427
$ws.Range("I8").Value = 0.9924269108841142
$ws.Range("J8").Value = 0.9924269108841142
$ws.Range("AQ8").Value = 0.9924269108841142
$ws.Range("AR8").Value = 0.9924269108841142
$ws.Range("I9").Value = 0.9899025478454855
$ws.Range("J9").Value = 0.9899025478454855
$ws.Range("AQ9").Value = 0.9899025478454855
$ws.Range("AR9").Value = 0.9899025478454855
$ws.Range("I10").Value = 0.9873781848068569
$ws.Range("J10").Value = 0.9873781848068569
$ws.Range("AQ10").Value = 0.9873781848068569
$ws.Range("AR10").Value = 0.9873781848068569
$ws.Range("I11").Value = 0.9848538217682282
$ws.Range("J11").Value = 0.9848538217682282
$ws.Range("AQ11").Value = 0.9848538217682282
$ws.Range("AR11").Value = 0.9848538217682282
$ws.Range("I12").Value = 0.9823294587295996
$ws.Range("J12").Value = 0.9823294587295996
$ws.Range("AQ12").Value = 0.9823294587295996
$ws.Range("AR12").Value = 0.9823294587295996
$ws.Range("I13").Value = 0.979805095690971
$ws.Range("J13").Value = 0.979805095690971
$ws.Range("AQ13").Value = 0.979805095690971
$ws.Range("AR13").Value = 0.979805095690971
$ws.Range("I14").Value = 0.9772807326523424
$ws.Range("J14").Value = 0.9772807326523424
$ws.Range("AQ14").Value = 0.9772807326523424
$ws.Range("AR14").Value = 0.9772807326523424
$ws.Range("I15").Value = 0.9747563696137138
$ws.Range("J15").Value = 0.9747563696137138
$ws.Range("AQ15").Value = 0.9747563696137138
$ws.Range("AR15").Value = 0.9747563696137138
$ws.Range("I16").Value = 0.9722320065750851
$ws.Range("J16").Value = 0.9722320065750851
$ws.Range("AQ16").Value = 0.9722320065750851
$ws.Range("AR16").Value = 0.9722320065750851
$ws.Range("I17").Value = 0.9697076435364566
$ws.Range("J17").Value = 0.9697076435364566
$ws.Range("AQ17").Value = 0.9697076435364566
$ws.Range("AR17").Value = 0.9697076435364566
$ws.Range("I18").Value = 0.9671832804978279
$ws.Range("J18").Value = 0.9671832804978279
$ws.Range("AQ18").Value = 0.9671832804978279
$ws.Range("AR18").Value = 0.9671832804978279
$ws.Range("I19").Value = 0.9646589174591993
$ws.Range("J19").Value = 0.9646589174591993
$ws.Range("AQ19").Value = 0.9646589174591993
$ws.Range("AR19").Value = 0.9646589174591993
$ws.Range("I20").Value = 0.9621345544205706
$ws.Range("J20").Value = 0.9621345544205706
$ws.Range("AQ20").Value = 0.9621345544205706
$ws.Range("AR20").Value = 0.9621345544205706
$ws.Range("I21").Value = 0.959610191381942
$ws.Range("J21").Value = 0.959610191381942
$ws.Range("AQ21").Value = 0.959610191381942
$ws.Range("AR21").Value = 0.959610191381942
$ws.Range("I22").Value = 0.9570858283433133
$ws.Range("J22").Value = 0.9570858283433133
$ws.Range("AQ22").Value = 0.9570858283433133
$ws.Range("AR22").Value = 0.9570858283433133
$ws.Range("I23").Value = 0.9545614653046848
$ws.Range("J23").Value = 0.9545614653046848
$ws.Range("AQ23").Value = 0.9545614653046848
$ws.Range("AR23").Value = 0.9545614653046848
$ws.Range("I24").Value = 0.9433955618175414
$ws.Range("J24").Value = 0.9433955618175414
$ws.Range("AQ24").Value = 0.9433955618175414
$ws.Range("AR24").Value = 0.9433955618175414
$ws.Range("I25").Value = 0.9322296583303981
$ws.Range("J25").Value = 0.9322296583303981
$ws.Range("AQ25").Value = 0.9322296583303981
$ws.Range("AR25").Value = 0.9322296583303981
$ws.Range("I26").Value = 0.9210637548432546
$ws.Range("J26").Value = 0.9210637548432546
$ws.Range("AQ26").Value = 0.9210637548432546
$ws.Range("AR26").Value = 0.9210637548432546
$ws.Range("I27").Value = 0.9098978513561113
$ws.Range("J27").Value = 0.9098978513561113
$ws.Range("AQ27").Value = 0.9098978513561113
$ws.Range("AR27").Value = 0.9098978513561113
$ws.Range("I28").Value = 0.8987319478689679
$ws.Range("J28").Value = 0.8987319478689679
$ws.Range("AQ28").Value = 0.8987319478689679
$ws.Range("AR28").Value = 0.8987319478689679
$ws.Range("I29").Value = 0.8865445579429376
$ws.Range("J29").Value = 0.8865445579429376
$ws.Range("AQ29").Value = 0.8865445579429376
$ws.Range("AR29").Value = 0.8865445579429376
$ws.Range("I30").Value = 0.8743571680169073
$ws.Range("J30").Value = 0.8743571680169073
$ws.Range("AQ30").Value = 0.8743571680169073
$ws.Range("AR30").Value = 0.8743571680169073
$ws.Range("I31").Value = 0.8621697780908771
$ws.Range("J31").Value = 0.8621697780908771
$ws.Range("AQ31").Value = 0.8621697780908771
$ws.Range("AR31").Value = 0.8621697780908771
$ws.Range("I32").Value = 0.8499823881648468
$ws.Range("J32").Value = 0.8499823881648468
$ws.Range("AQ32").Value = 0.8499823881648468
$ws.Range("AR32").Value = 0.8499823881648468
$ws.Range("I33").Value = 0.8377949982388165
$ws.Range("J33").Value = 0.8377949982388165
$ws.Range("AQ33").Value = 0.8377949982388165
$ws.Range("AR33").Value = 0.8377949982388165
$ws.Range("I34").Value = 0.8264353645649877
$ws.Range("J34").Value = 0.8264353645649877
$ws.Range("AQ34").Value = 0.8264353645649877
$ws.Range("AR34").Value = 0.8264353645649877
$ws.Range("I35").Value = 0.8150757308911589
$ws.Range("J35").Value = 0.8150757308911589
$ws.Range("AQ35").Value = 0.8150757308911589
$ws.Range("AR35").Value = 0.8150757308911589
$ws.Range("I36").Value = 0.8037160972173301
$ws.Range("J36").Value = 0.8037160972173301
$ws.Range("AQ36").Value = 0.8037160972173301
$ws.Range("AR36").Value = 0.8037160972173301
$ws.Range("I37").Value = 0.7923564635435013
$ws.Range("J37").Value = 0.7923564635435013
$ws.Range("AQ37").Value = 0.7923564635435013
$ws.Range("AR37").Value = 0.7923564635435013
$ws.Range("I38").Value = 0.7809968298696724
$ws.Range("J38").Value = 0.7809968298696724
$ws.Range("AQ38").Value = 0.7809968298696724
$ws.Range("AR38").Value = 0.7809968298696724
$ws.Range("I39").Value = 0.7696371961958437
$ws.Range("J39").Value = 0.7696371961958437
$ws.Range("AQ39").Value = 0.7696371961958437
$ws.Range("AR39").Value = 0.7696371961958437
$ws.Range("I40").Value = 0.7582775625220148
$ws.Range("J40").Value = 0.7582775625220148
$ws.Range("AQ40").Value = 0.7582775625220148
$ws.Range("AR40").Value = 0.7582775625220148
$ws.Range("I41").Value = 0.746917928848186
$ws.Range("J41").Value = 0.746917928848186
$ws.Range("AQ41").Value = 0.746917928848186
$ws.Range("AR41").Value = 0.746917928848186
$ws.Range("I42").Value = 0.7355582951743572
$ws.Range("J42").Value = 0.7355582951743572
$ws.Range("AQ42").Value = 0.7355582951743572
$ws.Range("AR42").Value = 0.7355582951743572
$ws.Range("I43").Value = 0.7241986615005284
$ws.Range("J43").Value = 0.7241986615005284
$ws.Range("AQ43").Value = 0.7241986615005284
$ws.Range("AR43").Value = 0.7241986615005284
